$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers already exist (A1:F1). Add data for new columns G (Actual
# Result) and H (Pass / Fail) for every test-case row, 2..14. ---

# Row 2 - TC-01
$ws.Range("G2").Value = "As Expected"
$ws.Range("H2").Value = "Pass"

# Row 3 - TC-02
$ws.Range("G3").Value = "Not working"
$ws.Range("H3").Value = "Fail"

# Row 4 - TC-03 (Actual Result styled with plain black Calibri font)
$ws.Range("G4").Value = "As Expected"
$ws.Range("G4").Font.Color = 0
$ws.Range("H4").Value = "Pass "

# Row 5 - TC-04
$ws.Range("G5").Value = "As Expected"
$ws.Range("G5").Font.Color = 0
$ws.Range("H5").Value = "Pass "

# Row 6 - TC-05
$ws.Range("G6").Value = "As Expected"
$ws.Range("G6").Font.Color = 0
$ws.Range("H6").Value = "Pass "

# Row 7 - TC-06
$ws.Range("G7").Value = "As Expected"
$ws.Range("G7").Font.Color = 0
$ws.Range("H7").Value = "Pass "

# Row 8 - TC-07
$ws.Range("G8").Value = "As Expected"
$ws.Range("G8").Font.Color = 0
$ws.Range("H8").Value = "Pass "

# Row 9 - TC-08 (Actual Result wraps like the other description columns)
$ws.Range("G9").Value = "Fogot password link is missing."
$ws.Range("G9").WrapText = $true
$ws.Range("H9").Value = "Fail"

# Row 10 - TC-09
$ws.Range("G10").Value = "Missing"
$ws.Range("H10").Value = "Fail"

# Row 11 - TC-10
$ws.Range("G11").Value = "Missing"
$ws.Range("H11").Value = "Fail"

# Row 12 - TC-11
$ws.Range("G12").Value = "Missing"
$ws.Range("H12").Value = "Fail"

# Row 13 - TC-12
$ws.Range("G13").Value = "Missing"
$ws.Range("H13").Value = "Fail"

# Row 14 - TC-13
$ws.Range("G14").Value = "As Expected"
$ws.Range("H14").Value = "Pass "

# Move the active selection to H14, matching where the author ended up
# after filling in the last row of results.
[void]$ws.Range("H14").Select()
